# Apply the edits described by the diff:
#  - handout/notes master "auto date" footer field text: 2020-01-20 -> 2020-02-10
#  - table-of-contents slide (#2): bump week numbers 2/2/3 -> 3/3/4
#  - slide titles on slides 3,4,5: "2주차 ..." -> "3주차 ..."
#  - slide title on slide 6: "3주차 계획" -> "4주차 계획"
#  - slide 5 picture ("그림 3") repositioned from (50800,-11296) to (0,0) EMU

$p = $ppt.ActivePresentation

# --- 1. Handout master footer date field -> 2020-02-10 --------------------
$handoutMaster = $p.HandoutMaster
$hmDateTime = $handoutMaster.HeadersFooters.DateAndTime
$hmDateTime.UseFormat = 0
$hmDateTime.Text = "2020-02-10"

# --- 2. Notes master footer date field -> 2020-02-10 -----------------------
$notesMaster = $p.NotesMaster
$nmDateTime = $notesMaster.HeadersFooters.DateAndTime
$nmDateTime.UseFormat = 0
$nmDateTime.Text = "2020-02-10"

# --- 3. Slide 2 : table of contents -----------------------------------------
$slide2 = $p.Slides.Item(2)
$tocShape = $slide2.Shapes.Item(2)
$tocText = $tocShape.TextFrame.TextRange
$tocText.Paragraphs(1).Characters(1, 1).Text = "3"   # 2주차 계획 -> 3주차 계획
$tocText.Paragraphs(2).Characters(1, 1).Text = "3"   # 2주차 현황 -> 3주차 현황
$tocText.Paragraphs(3).Characters(1, 1).Text = "4"   # 3주차 계획 -> 4주차 계획

# --- 4. Slide 3 title : 2주차 계획 -> 3주차 계획 -----------------------------
$slide3 = $p.Slides.Item(3)
$slide3.Shapes.Item(1).TextFrame.TextRange.Characters(1, 1).Text = "3"

# --- 5. Slide 4 title : 2주차 현황 -> 3주차 현황 -----------------------------
$slide4 = $p.Slides.Item(4)
$slide4.Shapes.Item(1).TextFrame.TextRange.Characters(1, 1).Text = "3"

# --- 6. Slide 5 title : 2주차 현황 -> 3주차 현황 -----------------------------
$slide5 = $p.Slides.Item(5)
$slide5.Shapes.Item(1).TextFrame.TextRange.Characters(1, 1).Text = "3"

# --- 7. Slide 5 picture ("그림 3") offset -> (0,0) --------------------------
$pic = $slide5.Shapes.Item(3)
$pic.Left = 0
$pic.Top = 0

# --- 8. Slide 6 title : 3주차 계획 -> 4주차 계획 -----------------------------
$slide6 = $p.Slides.Item(6)
$slide6.Shapes.Item(1).TextFrame.TextRange.Characters(1, 1).Text = "4"
